$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 3 (shifts existing rows 3..11 down to 4..12)
$ws.Rows.Item(3).Insert()

# Populate the newly inserted row 3 with "CHP" / "Point_2 Point_9"
$ws.Range("A3").Value = "CHP"
$ws.Range("B3").Value = "Point_2 Point_9"

# Update the active selection as recorded in the workbook
$ws.Range("C9").Select()
